$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641;  E = 0.4942365360607697; F = 1; G = 6.189590430959694 }
    3 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 0.7527432677738641;  E = 0.4942365360607697; F = 1; G = 4.358119930609447 }
    4 = @{ B = 0.6606524410359556; C = 1.655778082260271; D = 0.7527432677738641;  E = 0.4942365360607697; F = 1; G = 3.56341032713086 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 22.3905356188092;    E = 0.4942365360607697; F = 1; G = 27.82738278199502 }
    6 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 3.537761648806719;   E = 0.4942365360607697; F = 1; G = 8.974608811992548 }
    7 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 3.537761648806719;   E = 0.4942365360607697; F = 1; G = 7.143138311642302 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
